# Fixing the huge difference between the hourly and quarterly forecasted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Volue Forecast (MW)" values (column D) for rows 2-97
$forecastValues = @(
    1590.091708024305,
    1585.331061576131,
    1580.570415127957,
    1575.809768679783,
    1565.868886643684,
    1555.928004607585,
    1545.987122571486,
    1536.046240535386,
    1535.022124862358,
    1533.998009189331,
    1532.973893516303,
    1531.949777843275,
    1539.442725135224,
    1546.935672427172,
    1554.428619719121,
    1561.92156701107,
    1599.135743131193,
    1636.349919251316,
    1673.56409537144,
    1710.778271491563,
    1800.25217553843,
    1889.726079585296,
    1979.199983632163,
    2068.67388767903,
    2087.143764124873,
    2105.613640570715,
    2124.083517016557,
    2142.5533934624,
    2121.577913622901,
    2100.602433783401,
    2079.626953943902,
    2058.651474104402,
    2040.599198224548,
    2022.546922344693,
    2004.494646464838,
    1986.442370584983,
    1959.119362148208,
    1931.796353711433,
    1904.473345274658,
    1877.150336837884,
    1874.656664888733,
    1872.162992939583,
    1869.669320990433,
    1867.175649041283,
    1851.211375989515,
    1835.247102937748,
    1819.28282988598,
    1803.318556834212,
    1812.577357946448,
    1821.836159058684,
    1831.09496017092,
    1840.353761283156,
    1855.518627537862,
    1870.683493792568,
    1885.848360047275,
    1901.013226301981,
    1915.354822870217,
    1929.696419438452,
    1944.038016006688,
    1958.379612574924,
    1994.13815243752,
    2029.896692300116,
    2065.655232162712,
    2101.413772025308,
    2167.442387173591,
    2233.471002321874,
    2299.499617470156,
    2365.528232618439,
    2440.946894770997,
    2516.365556923554,
    2591.784219076112,
    2667.202881228669,
    2669.028392272934,
    2670.853903317199,
    2672.679414361464,
    2674.504925405729,
    2656.142431962624,
    2637.779938519519,
    2619.417445076414,
    2601.054951633309,
    2543.819811252643,
    2486.584670871978,
    2429.349530491313,
    2372.114390110648,
    2290.574896810305,
    2209.035403509962,
    2127.495910209619,
    2045.956416909276,
    2002.514026639272,
    1959.071636369269,
    1915.629246099265,
    1872.186855829261,
    1805.271238754834,
    1738.355621680406,
    1671.440004605979,
    1604.524387531552
)

$startRow = 2
for ($i = 0; $i -lt $forecastValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $forecastValues[$i]
}

# Previously-missing "Hydro Reservoir Actual (MW)", "Hydro River Actual (MW)"
# and "Hydro Actual (MW)" readings for rows 43-50 (columns B, C, E)
$actualValues = @(
    @{ Row = 43; B = 314; C = 1785; E = 2099 },
    @{ Row = 44; B = 313; C = 1784; E = 2097 },
    @{ Row = 45; B = 313; C = 1750; E = 2063 },
    @{ Row = 46; B = 216; C = 1697; E = 1913 },
    @{ Row = 47; B = 213; C = 1672; E = 1885 },
    @{ Row = 48; B = 212; C = 1664; E = 1876 },
    @{ Row = 49; B = 211; C = 1631; E = 1842 },
    @{ Row = 50; B = 224; C = 1649; E = 1873 }
)

foreach ($entry in $actualValues) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
    $ws.Cells.Item($entry.Row, 3).Value = $entry.C
    $ws.Cells.Item($entry.Row, 5).Value = $entry.E
}

